# Update the "Runmode" column (D) on the "Test Cases" sheet: flip every
# testcase to "N" except the two profile-specific ones (rows 12 and 15)
# that should keep running ("Y") — "running specific profile testcases".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$rowsToSkip = @(12, 15)

for ($r = 2; $r -le 20; $r++) {
    if ($rowsToSkip -contains $r) {
        continue
    }
    $ws.Cells.Item($r, 4).Value = "N"
}

# Update the sheet view's scroll position / selection to match the new
# working location.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B12").Select()
